# eCH-0160 zu EAD_v1.xlsx korrigieren
#
# Renames the (only) worksheet from "eCH-0160 zu xIsadg" to "eCH-0160 zu EAD".
# The workbook-level Print_Area / _xlnm.Print_Area defined names embed the
# sheet name in their reference string, so they are refreshed explicitly as
# well (re-asserting the unchanged $A$1:$J$42 area against the new sheet
# name). Finally the saved selection/scroll position is updated to match
# the author's last cursor position (cell C9, scrolled back to the top of
# the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "eCH-0160 zu EAD"

$ws.PageSetup.PrintArea = '$A$1:$J$42'

$ws.Range("C9").Select() | Out-Null
